# Auto-generated: Add data for 2025-07-19
# Updates 2025 (column L) violent-crime counts across Citywide Totals,
# By Neighborhood summary, and per-neighborhood sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 3639
$ws.Range("L3").Value = 3799
$ws.Range("L4").Value = 944
$ws.Range("L5").Value = 227
$ws.Range("L6").Value = 3321
$ws.Range("L7").Value = 11930

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("L6").Value = 46
$ws.Range("L7").Value = 134

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 221
$ws.Range("L3").Value = 257
$ws.Range("L6").Value = 211
$ws.Range("L7").Value = 767

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L2").Value = 154
$ws.Range("L5").Value = 11
$ws.Range("L7").Value = 561

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L2").Value = 66
$ws.Range("L3").Value = 53
$ws.Range("L7").Value = 169

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L2").Value = 131
$ws.Range("L5").Value = 15
$ws.Range("L7").Value = 429

$ws = $wb.Worksheets.Item("New City")
$ws.Range("L2").Value = 80
$ws.Range("L6").Value = 66
$ws.Range("L7").Value = 226

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L3").Value = 83
$ws.Range("L7").Value = 200

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L7").Value = 400
$ws.Range("L8").Value = 767
$ws.Range("L10").Value = 75
$ws.Range("L11").Value = 200
$ws.Range("L14").Value = 60
$ws.Range("L20").Value = 300
$ws.Range("L23").Value = 130
$ws.Range("L25").Value = 61
$ws.Range("L27").Value = 109
$ws.Range("L29").Value = 651
$ws.Range("L32").Value = 14
$ws.Range("L33").Value = 561
$ws.Range("L37").Value = 429
$ws.Range("L46").Value = 28
$ws.Range("L49").Value = 63
$ws.Range("L51").Value = 147
$ws.Range("L53").Value = 134
$ws.Range("L56").Value = 9
$ws.Range("L63").Value = 39
$ws.Range("L64").Value = 77
$ws.Range("L65").Value = 226
$ws.Range("L67").Value = 426
$ws.Range("L73").Value = 102
$ws.Range("L76").Value = 178
$ws.Range("L78").Value = 151
$ws.Range("L79").Value = 315
$ws.Range("L85").Value = 617
$ws.Range("L91").Value = 171
$ws.Range("L93").Value = 64
$ws.Range("L94").Value = 144
$ws.Range("L95").Value = 169
$ws.Range("L99").Value = 200
$ws.Range("L101").Value = 11930

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L3").Value = 162
$ws.Range("L6").Value = 99
$ws.Range("L7").Value = 426

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("L4").Value = 7
$ws.Range("L6").Value = 27
$ws.Range("L7").Value = 63

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 198
$ws.Range("L3").Value = 245
$ws.Range("L5").Value = 10
$ws.Range("L6").Value = 165
$ws.Range("L7").Value = 651

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L2").Value = 35
$ws.Range("L3").Value = 32
$ws.Range("L7").Value = 178

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("L6").Value = 14
$ws.Range("L7").Value = 60

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("L2").Value = 33
$ws.Range("L7").Value = 75

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L2").Value = 43
$ws.Range("L4").Value = 17
$ws.Range("L7").Value = 151

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("L2").Value = 8
$ws.Range("L7").Value = 28

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("L3").Value = 51
$ws.Range("L7").Value = 130

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("L3").Value = 73
$ws.Range("L7").Value = 171

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L2").Value = 108
$ws.Range("L7").Value = 315

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("L6").Value = 24
$ws.Range("L7").Value = 77

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L4").Value = 28
$ws.Range("L7").Value = 300

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("L2").Value = 22
$ws.Range("L7").Value = 64

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L2").Value = 131
$ws.Range("L6").Value = 112
$ws.Range("L7").Value = 400

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L2").Value = 35
$ws.Range("L3").Value = 33
$ws.Range("L7").Value = 144

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("L3").Value = 30
$ws.Range("L7").Value = 61

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L4").Value = 15
$ws.Range("L6").Value = 49
$ws.Range("L7").Value = 200

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("L3").Value = 30
$ws.Range("L7").Value = 102

$ws = $wb.Worksheets.Item("Galewood")
$ws.Range("L2").Value = 9
$ws.Range("L7").Value = 14

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("L6").Value = 32
$ws.Range("L7").Value = 109

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L3").Value = 45
$ws.Range("L7").Value = 147

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L2").Value = 180
$ws.Range("L3").Value = 255
$ws.Range("L4").Value = 42
$ws.Range("L5").Value = 12
$ws.Range("L6").Value = 128
$ws.Range("L7").Value = 617

$ws = $wb.Worksheets.Item("Magnificent Mile")
$ws.Range("L6").Value = 6
$ws.Range("L7").Value = 9
